# Clean up of XLPNR
# ------------------
# The "Gearbox Tests" sheet had three mass measurements (HOUSING.mass,
# FASTENERS.mass, GEARS.mass) recorded in C5:C7 that were off by a factor
# of 1000 (grams instead of kilograms, per the "kg" unit label in column D).
# Correct each value by shifting the decimal point (dividing by 1000),
# preserving full floating point precision.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gearbox Tests")

$ws.Range("C5").Value = 128.85048535948661   # HOUSING.mass,   was 128850.48535948661
$ws.Range("C6").Value = 3.8029230326870378   # FASTENERS.mass, was 3802.9230326870379
$ws.Range("C7").Value = 54.455593060061851   # GEARS.mass,     was 54455.593060061852
